# Functions for button which opens input or output files
#
# 1. Fix a typo in the shared-string account number on the "input" sheet:
#    "77x1240 1268 1111 0010 6850 0503" -> "77 1240 1268 1111 0010 6850 0503"
# 2. Move the sheet's active cell/selection from A4 to A13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# --- Fix the mistyped account number (the "x" should be a space) ---
$ws.Range("A2").Value = "77 1240 1268 1111 0010 6850 0503"

# --- Update the saved selection on the sheet ---
$ws.Activate()
$ws.Range("A13").Select()
